# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly observation (week of 2022-01-18, serial 44579) is inserted
# for the "Betarraga" subset right before the existing row that used to sit
# at row 94 (date serial 44175). Inserting 2 new rows at row 94 pushes every
# subsequent record down by one weekly pair (2 rows); the last existing pair
# (originally rows 226/227, date serial 44511) ends up re-appearing as the
# new final pair at rows 228/229 — i.e. nothing is deleted, the whole table
# just grows by one row-pair at the top of this date-ordered run and the
# overall sheet dimension grows from R227 to R229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at 94/95 - everything currently at/after row 94
# (including the old row 94/95 pair) shifts down to 96/97, ... and the old
# 226/227 pair ends up at 228/229.
$ws.Rows.Item(94).Resize(2).Insert()

# Row 94 - "Primera" quality
$ws.Cells.Item(94, 1).Value = 8
$ws.Cells.Item(94, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(94, 3).Value = "Coquimbo"
$ws.Cells.Item(94, 4).Value = 44579
$ws.Cells.Item(94, 5).Value = 4
$ws.Cells.Item(94, 6).Value = 100114014
$ws.Cells.Item(94, 7).Value = "Betarraga"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 3120
$ws.Cells.Item(94, 11).Value = 450
$ws.Cells.Item(94, 12).Value = 500
$ws.Cells.Item(94, 13).Value = 475
$ws.Cells.Item(94, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 158
$ws.Cells.Item(94, 17).Value = 3
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# Row 95 - "Segunda" quality
$ws.Cells.Item(95, 1).Value = 8
$ws.Cells.Item(95, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44579
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = 100114014
$ws.Cells.Item(95, 7).Value = "Betarraga"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Segunda"
$ws.Cells.Item(95, 10).Value = 1560
$ws.Cells.Item(95, 11).Value = 350
$ws.Cells.Item(95, 12).Value = 400
$ws.Cells.Item(95, 13).Value = 375
$ws.Cells.Item(95, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(95, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(95, 16).Value = 125
$ws.Cells.Item(95, 17).Value = 3
$ws.Cells.Item(95, 18).Value = "Hortaliza"

# Make sure the date cells carry the same date-time number format the rest
# of column D uses (style index 2 in the original file), matching every
# other row in this column.
$ws.Range("D94:D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
